$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The commit reorders the four data records currently sitting in rows 2-5
# (row 1 is the header, row 6 is untouched): the record in row 2 swaps
# completely with the one in row 5, and the record in row 3 swaps
# completely with the one in row 4. Row 3 <-> row 4 share an identical set
# of populated columns, but row 2 and row 5 differ - row 2 currently has a
# few extra populated cells (J/K/L/N/AC/AF) that row 5 doesn't have, so
# those need to move across too. We therefore write each changed cell
# explicitly instead of blitting whole-row ranges, so that:
#   - cells that stay blank on both sides (I2:I5, AT2:AT5) are left alone
#   - cells that need to disappear are actually cleared (not just blanked)
#   - cells that need to newly appear (but stay blank) are (re)created
#   - text-like cells (dates in Y/AA) are never touched, avoiding any
#     auto type coercion

function Swap-Cell($addr1, $addr2) {
    $v1 = $ws.Range($addr1).Value()
    $v2 = $ws.Range($addr2).Value()
    $ws.Range($addr1).Value = $v2
    $ws.Range($addr2).Value = $v1
}

# --- Row 2 <-> Row 5 -------------------------------------------------
foreach ($col in @("A","B","D","E","F","G","H","Q","R")) {
    Swap-Cell "$col`2" "$col`5"
}

# Row 2 loses J2,K2,L2,N2,AC2,AF2 (they belong to the record moving to row 5)
foreach ($col in @("J","K","L","N","AC","AF")) {
    $ws.Range("$col`2").ClearContents()
}

# Row 5 gains the same columns back. K/AC carry real text, the rest stay
# blank but present (mirrored via a style touch so the cell persists as an
# empty cell rather than being dropped).
$ws.Range("K5").Value = "blomning"
$ws.Range("AC5").Value = "Blommande"
foreach ($col in @("J","L","N","AF")) {
    $ws.Range("$col`5").Style = "Normal"
}

# --- Row 3 <-> Row 4 (identical column layout, straight swap) --------
foreach ($col in @("A","B","D","E","F","G","H","Q","R")) {
    Swap-Cell "$col`3" "$col`4"
}

Write-Host "done"
